# msz - mandatory fields checks part 1
#
# Adds 3 new rows to the "Tabelle1" test-step sheet describing a new
# mandatory-fields check flow for the vehicle data page, and moves the
# illustrative screenshot picture down to make room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- New test rows -----------------------------------------------------

$ws.Range("A17").Value = "Vehicle Page check for open mandatory fields"
$ws.Range("B17").Value = "<CHK>"
$ws.Range("C17").Value = "Vehicle Page check for open mandatory fields"
$ws.Range("H17").Value = "<NOP>"

$ws.Range("A18").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_FillMake"
$ws.Range("B18").Value = "<SET>"
$ws.Range("C18").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_FillMake"
$ws.Range("H18").Value = "<NOP>"

$ws.Range("A19").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_CheckFilledMake"
$ws.Range("B19").Value = "<CHK>"
$ws.Range("C19").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields_CheckFilledMake"
$ws.Range("H19").Value = "<NOP>"

# --- Column widths (auto-fit growth caused by the new, longer text) ----

$ws.Columns("A:A").ColumnWidth = 70.6
$ws.Columns("C:C").ColumnWidth = 70.6

# --- Move the screenshot picture down so it still sits below the table -

$shp = $ws.Shapes.Item(1)
$shp.Top = 337.80007874015746
$shp.Left = 0
$shp.Width = 676.6974015748032
$shp.Height = 398.4

# --- Selection left where the author ended up editing ------------------

$null = $ws.Range("D22").Select()
